$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 223.28572
$ws.Range("I5").Value = 354
$ws.Range("J5").Value = 49
$ws.Range("K5").Value = 354
$ws.Range("L5").Value = 49
$ws.Range("M5").Value = -239
$ws.Range("N5").Value = -279
$ws.Range("H6").Value = 143259.14
$ws.Range("I6").Value = 200062.8
$ws.Range("K6").Value = 600188.3999999999
$ws.Range("M6").Value = -600076.3999999999
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 1000
$ws.Range("K29").Value = 3000
$ws.Range("M29").Value = -2719
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H38").Value = 341.33334
$ws.Range("I38").Value = 281.45456
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 844.36368
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -472.36368
$ws.Range("N38").Value = -3744
$ws.Range("H70").Value = 9000
$ws.Range("J70").Value = 7000
$ws.Range("L70").Value = 21000
$ws.Range("N70").Value = -21540
$ws.Range("H73").Value = 9000
$ws.Range("J73").Value = 7000
$ws.Range("L73").Value = 21000
$ws.Range("N73").Value = -22872
$ws.Range("H74").Value = 8477.666999999999
$ws.Range("I74").Value = 8477.666999999999
$ws.Range("K74").Value = 8477.666999999999
$ws.Range("M74").Value = -7541.666999999999
$ws.Range("H77").Value = 8477.666999999999
$ws.Range("I77").Value = 8477.666999999999
$ws.Range("K77").Value = 42388.335
$ws.Range("M77").Value = -37708.335
$ws.Range("H86").Value = 1224.5
$ws.Range("I86").Value = 1132.6666
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 1132.6666
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -9.666600000000017
$ws.Range("N86").Value = -3746
$ws.Range("H88").Value = 2126.125
$ws.Range("I88").Value = 2648.75
$ws.Range("J88").Value = 1603.5
$ws.Range("K88").Value = 2648.75
$ws.Range("L88").Value = 1603.5
$ws.Range("M88").Value = -2242.75
$ws.Range("N88").Value = -2415.5
$ws.Range("H89").Value = 1224.5
$ws.Range("I89").Value = 1132.6666
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 5663.333000000001
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = -47.33300000000054
$ws.Range("N89").Value = -18732
$ws.Range("H91").Value = 2126.125
$ws.Range("I91").Value = 2648.75
$ws.Range("J91").Value = 1603.5
$ws.Range("K91").Value = 2648.75
$ws.Range("L91").Value = 1603.5
$ws.Range("M91").Value = -1244.75
$ws.Range("N91").Value = -4411.5
$ws.Range("H98").Value = 1363.9333
$ws.Range("I98").Value = 963.6667
$ws.Range("J98").Value = 2965
$ws.Range("K98").Value = 963.6667
$ws.Range("L98").Value = 2965
$ws.Range("M98").Value = 534.3333
$ws.Range("N98").Value = -5961
$ws.Range("H113").Value = 4763.5454
$ws.Range("I113").Value = 4990
$ws.Range("J113").Value = 2499
$ws.Range("K113").Value = 4990
$ws.Range("L113").Value = 2499
$ws.Range("M113").Value = -1736
$ws.Range("N113").Value = -9007
$ws.Range("H116").Value = 2300.1667
$ws.Range("I116").Value = 2440.2
$ws.Range("J116").Value = 1600
$ws.Range("K116").Value = 2440.2
$ws.Range("L116").Value = 1600
$ws.Range("M116").Value = 1001.8
$ws.Range("N116").Value = -8484
$ws.Range("H122").Value = 1363.9333
$ws.Range("I122").Value = 963.6667
$ws.Range("J122").Value = 2965
$ws.Range("K122").Value = 2891.0001
$ws.Range("L122").Value = 8895
$ws.Range("M122").Value = -441.0001000000002
$ws.Range("N122").Value = -13795
$ws.Range("H127").Value = 1783.5
$ws.Range("I127").Value = 1783.5
$ws.Range("K127").Value = 5350.5
$ws.Range("M127").Value = -390.5
$ws.Range("H129").Value = 1534.3529
$ws.Range("I129").Value = 722.6667
$ws.Range("J129").Value = 2447.5
$ws.Range("K129").Value = 2168.0001
$ws.Range("L129").Value = 7342.5
$ws.Range("M129").Value = 2831.9999
$ws.Range("N129").Value = -17342.5
$ws.Range("H131").Value = 948.44446
$ws.Range("I131").Value = 1040.75
$ws.Range("J131").Value = 210
$ws.Range("K131").Value = 3122.25
$ws.Range("L131").Value = 630
$ws.Range("M131").Value = 1917.75
$ws.Range("N131").Value = -10710
$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120
$ws.Range("H136").Value = 199663.33
$ws.Range("J136").Value = 199663.33
$ws.Range("L136").Value = 199663.33
$ws.Range("N136").Value = -209863.33
$ws.Range("H137").Value = 12989.714
$ws.Range("I137").Value = 9379.4
$ws.Range("K137").Value = 28138.2
$ws.Range("M137").Value = -25588.2
$ws.Range("H138").Value = 2322.5
$ws.Range("J138").Value = 4250
$ws.Range("L138").Value = 12750
$ws.Range("N138").Value = -23030
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2604.0344
$ws.Range("I32").Value = 2604.0344
$ws.Range("K32").Value = 2604.0344
$ws.Range("M32").Value = -2317.0344
$ws.Range("H33").Value = 1500
$ws.Range("I33").Value = 1500
$ws.Range("K33").Value = 1500
$ws.Range("M33").Value = -1171
$ws.Range("H61").Value = 5119.8
$ws.Range("J61").Value = 8501.833000000001
$ws.Range("L61").Value = 8501.833000000001
$ws.Range("N61").Value = -8925.833000000001
$ws.Range("H88").Value = 252722.75
$ws.Range("I88").Value = 502005.5
$ws.Range("J88").Value = 3440
$ws.Range("K88").Value = 502005.5
$ws.Range("L88").Value = 3440
$ws.Range("M88").Value = -501599.5
$ws.Range("N88").Value = -4252
$ws.Range("H91").Value = 252722.75
$ws.Range("I91").Value = 502005.5
$ws.Range("J91").Value = 3440
$ws.Range("K91").Value = 502005.5
$ws.Range("L91").Value = 3440
$ws.Range("M91").Value = -500601.5
$ws.Range("N91").Value = -6248
$ws.Range("H110").Value = 7879.8423
$ws.Range("I110").Value = 8519.462
$ws.Range("K110").Value = 8519.462
$ws.Range("M110").Value = -6474.462
$ws.Range("H122").Value = 1140.1333
$ws.Range("I122").Value = 1140.1333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3420.3999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -970.3998999999999
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 3255.8838
$ws.Range("I132").Value = 3214.1052
$ws.Range("J132").Value = 3573.4
$ws.Range("K132").Value = 9642.3156
$ws.Range("L132").Value = 10720.2
$ws.Range("M132").Value = -7112.3156
$ws.Range("N132").Value = -15780.2
$ws.Range("H136").Value = 5119.8
$ws.Range("J136").Value = 8501.833000000001
$ws.Range("L136").Value = 25505.499
$ws.Range("N136").Value = -30605.499
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6577.143
$ws.Range("I86").Value = 3453.5557
$ws.Range("J86").Value = 12199.6
$ws.Range("K86").Value = 3453.5557
$ws.Range("L86").Value = 12199.6
$ws.Range("M86").Value = -2330.5557
$ws.Range("N86").Value = -14445.6
$ws.Range("H89").Value = 6577.143
$ws.Range("I89").Value = 3453.5557
$ws.Range("J89").Value = 12199.6
$ws.Range("K89").Value = 17267.7785
$ws.Range("L89").Value = 60998
$ws.Range("M89").Value = -11651.7785
$ws.Range("N89").Value = -72230
$ws.Range("H107").Value = 1145.6
$ws.Range("I107").Value = 707.25
$ws.Range("K107").Value = 707.25
$ws.Range("M107").Value = 1212.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 223.54546
$ws.Range("I7").Value = 140
$ws.Range("K7").Value = 140
$ws.Range("M7").Value = -27
$ws.Range("H16").Value = 2401.8572
$ws.Range("I16").Value = 1160
$ws.Range("K16").Value = 1160
$ws.Range("M16").Value = -873
$ws.Range("H31").Value = 1740.7142
$ws.Range("I31").Value = 1740.7142
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1740.7142
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1445.7142
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 1740.7142
$ws.Range("I34").Value = 1740.7142
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1740.7142
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1538.7142
$ws.Range("N34").ClearContents()
$ws.Range("H58").Value = 2197.0625
$ws.Range("I58").Value = 1622.9231
$ws.Range("J58").Value = 4685
$ws.Range("K58").Value = 1622.9231
$ws.Range("L58").Value = 4685
$ws.Range("M58").Value = -1419.9231
$ws.Range("N58").Value = -5091
$ws.Range("H62").Value = 6450
$ws.Range("I62").Value = 6450
$ws.Range("K62").Value = 6450
$ws.Range("M62").Value = -5826
$ws.Range("H65").Value = 6450
$ws.Range("I65").Value = 6450
$ws.Range("K65").Value = 32250
$ws.Range("M65").Value = -29130
$ws.Range("H99").Value = 4123
$ws.Range("I99").Value = 3544.2
$ws.Range("K99").Value = 3544.2
$ws.Range("M99").Value = -2046.2
$ws.Range("H113").Value = 2401.8572
$ws.Range("I113").Value = 1160
$ws.Range("K113").Value = 1160
$ws.Range("M113").Value = 1010
$ws.Range("H126").Value = 4123
$ws.Range("I126").Value = 3544.2
$ws.Range("K126").Value = 10632.6
$ws.Range("M126").Value = -8162.599999999999
$ws.Range("H132").Value = 743.36365
$ws.Range("I132").Value = 743.36365
$ws.Range("K132").Value = 2230.09095
$ws.Range("M132").Value = 299.9090500000002
$ws.Range("H134").Value = 2555.7144
$ws.Range("I134").Value = 2616.7693
$ws.Range("J134").Value = 1762
$ws.Range("K134").Value = 7850.3079
$ws.Range("L134").Value = 5286
$ws.Range("M134").Value = -5315.3079
$ws.Range("N134").Value = -10356
$ws.Range("H136").Value = 2197.0625
$ws.Range("I136").Value = 1622.9231
$ws.Range("J136").Value = 4685
$ws.Range("K136").Value = 4868.7693
$ws.Range("L136").Value = 14055
$ws.Range("M136").Value = -2318.7693
$ws.Range("N136").Value = -19155
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 93720.73
$ws.Range("I7").Value = 167016.67
$ws.Range("J7").Value = 5765.6
$ws.Range("K7").Value = 501050.01
$ws.Range("L7").Value = 17296.8
$ws.Range("M7").Value = -500938.01
$ws.Range("N7").Value = -17520.8
$ws.Range("H14").Value = 483.8889
$ws.Range("I14").Value = 483.8889
$ws.Range("K14").Value = 1451.6667
$ws.Range("M14").Value = -1278.6667
$ws.Range("H55").Value = 15832.167
$ws.Range("J55").Value = 15832.167
$ws.Range("L55").Value = 47496.501
$ws.Range("N55").Value = -47850.501
$ws.Range("H92").Value = 154.77777
$ws.Range("I92").Value = 170.8
$ws.Range("J92").Value = 134.75
$ws.Range("K92").Value = 512.4000000000001
$ws.Range("L92").Value = 404.25
$ws.Range("M92").Value = 735.5999999999999
$ws.Range("N92").Value = -2900.25
$ws.Range("H132").Value = 1433.1666
$ws.Range("J132").Value = 1450
$ws.Range("L132").Value = 13050
$ws.Range("N132").Value = -18110
$ws.Range("H139").Value = 5167.5293
$ws.Range("I139").Value = 5121
$ws.Range("K139").Value = 15363
$ws.Range("M139").Value = -10223
$ws.Range("H140").Value = 558189.1
$ws.Range("I140").Value = 558189.1
$ws.Range("K140").Value = 1674567.3
$ws.Range("M140").Value = -1669387.3
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H97").Value = 575.6923
$ws.Range("I97").Value = 598.6667
$ws.Range("K97").Value = 598.6667
$ws.Range("M97").Value = -102.6667
$ws.Range("H107").Value = 2940.3635
$ws.Range("I107").Value = 1792.4286
$ws.Range("J107").Value = 4949.25
$ws.Range("K107").Value = 1792.4286
$ws.Range("L107").Value = 4949.25
$ws.Range("M107").Value = 127.5714
$ws.Range("N107").Value = -8789.25
$ws.Range("H113").Value = 1654.75
$ws.Range("I113").Value = 1651.3334
$ws.Range("K113").Value = 1651.3334
$ws.Range("M113").Value = 518.6666
$ws.Range("H122").Value = 3438.4119
$ws.Range("I122").Value = 3027.1538
$ws.Range("K122").Value = 9081.4614
$ws.Range("M122").Value = -6631.4614
$ws.Range("H126").Value = 4813.3335
$ws.Range("J126").Value = 4998
$ws.Range("L126").Value = 14994
$ws.Range("N126").Value = -19934
$ws.Range("H132").Value = 2062.8667
$ws.Range("I132").Value = 2067.3572
$ws.Range("K132").Value = 6202.071599999999
$ws.Range("M132").Value = -3672.071599999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4474.7
$ws.Range("I7").Value = 4068.5
$ws.Range("K7").Value = 4068.5
$ws.Range("M7").Value = -3956.5
$ws.Range("H16").Value = 531
$ws.Range("I16").Value = 483.75
$ws.Range("J16").Value = 720
$ws.Range("K16").Value = 483.75
$ws.Range("L16").Value = 720
$ws.Range("M16").Value = -313.75
$ws.Range("N16").Value = -1060
$ws.Range("H40").Value = 2422.6
$ws.Range("I40").Value = 1802.8889
$ws.Range("K40").Value = 1802.8889
$ws.Range("M40").Value = -1666.8889
$ws.Range("H55").Value = 345.75
$ws.Range("I55").Value = 351
$ws.Range("J55").Value = 330
$ws.Range("K55").Value = 351
$ws.Range("L55").Value = 330
$ws.Range("M55").Value = -178
$ws.Range("N55").Value = -676
$ws.Range("H61").Value = 1611.75
$ws.Range("I61").Value = 1599.3334
$ws.Range("K61").Value = 1599.3334
$ws.Range("M61").Value = -1397.3334
$ws.Range("H93").Value = 1576.625
$ws.Range("I93").Value = 1785.5
$ws.Range("J93").Value = 950
$ws.Range("K93").Value = 1785.5
$ws.Range("L93").Value = 950
$ws.Range("M93").Value = -537.5
$ws.Range("N93").Value = -3446
$ws.Range("H113").Value = 1611.75
$ws.Range("I113").Value = 1599.3334
$ws.Range("K113").Value = 1599.3334
$ws.Range("M113").Value = 570.6666
$ws.Range("H122").Value = 3373.4443
$ws.Range("I122").Value = 3195
$ws.Range("J122").Value = 3998
$ws.Range("K122").Value = 9585
$ws.Range("L122").Value = 11994
$ws.Range("M122").Value = -7135
$ws.Range("N122").Value = -16894
$ws.Range("H126").Value = 4474.7
$ws.Range("I126").Value = 4068.5
$ws.Range("K126").Value = 12205.5
$ws.Range("M126").Value = -9735.5
$ws.Range("H132").Value = 2256.257
$ws.Range("I132").Value = 2211.8696
$ws.Range("J132").Value = 2341.3333
$ws.Range("K132").Value = 6635.6088
$ws.Range("L132").Value = 7023.999899999999
$ws.Range("M132").Value = -4105.6088
$ws.Range("N132").Value = -12083.9999
$ws.Range("H136").Value = 25004616
$ws.Range("I136").Value = 4018
$ws.Range("K136").Value = 12054
$ws.Range("M136").Value = -9504
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 42455.5
$ws.Range("I45").Value = 31812.8
$ws.Range("J45").Value = 53098.2
$ws.Range("K45").Value = 31812.8
$ws.Range("L45").Value = 53098.2
$ws.Range("M45").Value = -31321.8
$ws.Range("N45").Value = -54080.2
$ws.Range("H113").Value = 918.4545000000001
$ws.Range("I113").Value = 567
$ws.Range("K113").Value = 1701
$ws.Range("M113").Value = 469
$ws.Range("H122").Value = 9072.75
$ws.Range("I122").Value = 8368.143
$ws.Range("K122").Value = 25104.429
$ws.Range("M122").Value = -22654.429
$ws.Range("H126").Value = 1739.1
$ws.Range("I126").Value = 1154.5555
$ws.Range("K126").Value = 3463.6665
$ws.Range("M126").Value = -993.6664999999998
$ws.Range("H132").Value = 19477
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 7780.9165
$ws.Range("I136").Value = 7780.9165
$ws.Range("K136").Value = 23342.7495
$ws.Range("M136").Value = -20792.7495
